$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two obsolete pupil rows (old rows 3 & 4) ---------------
$ws.Rows("3:4").Delete()

# --- Extend the table with two new columns: ID_LEHRER / ID_KATEGORIE ---
$ws.Range("E1").Value = "ID_LEHRER"
$ws.Range("F1").Value = "ID_KATEGORIE"

# --- Fill in the remaining values for the single pupil row left --------
$ws.Range("D2").Value = "FIAE17A"
$ws.Range("E2").Value = "HK"
$ws.Range("F2").Value = 17

# --- Birthdate format now shows a 4-digit year --------------------------
$ws.Range("A2").NumberFormat = "dd/mm/yyyy;@"

# --- Column widths -------------------------------------------------------
$ws.Columns("A").ColumnWidth = 9.3333333333333
$ws.Columns("B:C").ColumnWidth = 7.5
$ws.Columns("F").ColumnWidth = 13.3333333333333

# --- Selection matches the author's last active cell ---------------------
$null = $ws.Range("E2").Select()
